# Updated bus voltage magnitude (vm_pu) results for the 380 kV case.
# Slack/reference bus voltage setpoint changed from 1.05 pu to 1.02 pu,
# which ripples through the load-flow solution for every bus (columns
# B-F, I-M) across all rows (buses 0-23, sheet rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B","C","D","E","F","I","J","K","L","M")

$data = @{
    2 = @(1.02, 1.019430475293696, 1.024624822111795, 1.020565154313159, 1.030511454988552, 1.028794099063778, 1.024633284677864, 1.027452775983998, 1.023405064274411, 1.033322254128945)
    3 = @(1.02, 1.020374550968546, 1.025290576997821, 1.021365943612955, 1.031670877349935, 1.028939656907195, 1.025214140346149, 1.02792604620203, 1.024012134932392, 1.03428910392326)
    4 = @(1.02, 1.020985824918355, 1.025721490802771, 1.021884833328112, 1.032421622644335, 1.029032540051098, 1.025589802479125, 1.028231729751108, 1.024405025861916, 1.034914670501949)
    5 = @(1.02, 1.021242898067524, 1.025902675849082, 1.022103147084658, 1.032737360771329, 1.029071275866011, 1.025747684513912, 1.028360105661419, 1.024570214453358, 1.035177646557867)
    6 = @(1.02, 1.021286067247502, 1.025933099272853, 1.02213981300784, 1.032790381909693, 1.029077761456165, 1.025774190863559, 1.028381652691359, 1.02459795133099, 1.035221800677847)
    7 = @(1.02, 1.020989259575529, 1.025723911694316, 1.021887749770788, 1.032425841062289, 1.029033058868458, 1.025591912289586, 1.028233445641745, 1.024407233053279, 1.034918184448793)
    8 = @(1.02, 1.019749448109126, 1.02484979027708, 1.020835633889305, 1.030903179197614, 1.028843560451348, 1.024829626520697, 1.027612834058649, 1.023610210471685, 1.033649015758506)
    9 = @(1.02, 1.01756778961511, 1.023310495848176, 1.018987278934483, 1.028224050016764, 1.028499683461504, 1.02348495743812, 1.026515031778076, 1.022206371813244, 1.031412204358197)
    10 = @(1.02, 1.01611544009761, 1.022285058361629, 1.017758880743581, 1.026440649275957, 1.02826376325615, 1.022587598009859, 1.02578038529293, 1.021270944245854, 1.029920753611832)
    11 = @(1.02, 1.015487060264338, 1.021841227850909, 1.017227895276424, 1.025669052783165, 1.02816003061984, 1.022198823809673, 1.025461626382033, 1.02086601430635, 1.029274882203541)
    12 = @(1.02, 1.015253727268612, 1.02167639928823, 1.017030802814988, 1.025382541590321, 1.028121263062008, 1.022054384643775, 1.025343127929341, 1.020715623571476, 1.029034967245629)
    13 = @(1.02, 1.015303774573711, 1.021711754204382, 1.017073073500944, 1.025443994959025, 1.028129589536972, 1.022085368732402, 1.025368550635779, 1.020747882054983, 1.029086430251622)
    14 = @(1.02, 1.015467771347671, 1.021827602454617, 1.017211600706394, 1.025645367760227, 1.028156830906541, 1.022186885054023, 1.025451833241345, 1.020853582583931, 1.029255050953477)
    15 = @(1.02, 1.015568825161567, 1.02189898437645, 1.017296970365817, 1.025769452651867, 1.028173583868055, 1.022249428507162, 1.025503133546488, 1.020918710613528, 1.029358942465372)
    16 = @(1.02, 1.016157154119393, 1.022314518051136, 1.017794139958115, 1.026491870802102, 1.02827061444341, 1.022613395271999, 1.025801526598516, 1.021297820631635, 1.029963616631343)
    17 = @(1.02, 1.016526330818536, 1.022575223248984, 1.018106248146478, 1.026945192544524, 1.028331057053579, 1.022841645823865, 1.025988526629325, 1.021535658044133, 1.030342895762252)
    18 = @(1.02, 1.016741713264883, 1.02272730649828, 1.018288384173789, 1.027209668035006, 1.028366159961275, 1.022974760042715, 1.026097537609986, 1.021674395757967, 1.03056411674805)
    19 = @(1.02, 1.016815161224492, 1.022779166026578, 1.018350502846983, 1.027299857626078, 1.028378103310754, 1.023020145059238, 1.02613469683229, 1.02172170360449, 1.030639546351381)
    20 = @(1.02, 1.016486716691413, 1.022547250145098, 1.018072752719142, 1.026896549146336, 1.028324587878082, 1.022817158811775, 1.025968469803783, 1.021510139179011, 1.030302203333547)
    21 = @(1.02, 1.015419476303302, 1.021793487198706, 1.017170804046079, 1.025586065893477, 1.028148815531463, 1.022156991881077, 1.025427311266754, 1.020822455901866, 1.029205396623692)
    22 = @(1.02, 1.014748895588265, 1.021319739816483, 1.016604518934476, 1.024762657387571, 1.028036931556465, 1.021741739028541, 1.025086501427516, 1.020390188251336, 1.028515735390592)
    23 = @(1.02, 1.01510434162464, 1.021570865440124, 1.016904640635637, 1.025199110292002, 1.028096372998592, 1.021961889229142, 1.025267224182939, 1.020619331188763, 1.028881343113813)
    24 = @(1.02, 1.016504616462199, 1.022559889930168, 1.018087887587293, 1.026918528816746, 1.028327511490457, 1.022828223501777, 1.025977532814541, 1.021521670021237, 1.030320590509234)
    25 = @(1.02, 1.018131435311792, 1.023708311985212, 1.019464451429344, 1.028916195752301, 1.028589761200231, 1.023832750847977, 1.026799333271335, 1.022569219482204, 1.031990517198356)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}

Write-Output "Applied all vm_pu updates (380 kV case)"
